$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.838.02'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").Value = '3.818.23'
$ws.Range("E3").Value = '  +0.99%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.94'
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.18'
$ws.Range("E6").Value = '  -0.45%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.11'
$ws.Range("E13").Value = '  +0.27%  '

$ws.Range("D14").Value = '4.459.59'
$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = '3.824.01'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("D16").Value = '67.888.25'
$ws.Range("E16").Value = '  +0.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.38'
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.11'
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '463.83'
$ws.Range("E20").Value = '  +1.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.87'
$ws.Range("E21").Value = '  -1.47%  '

$ws.Range("E22").Value = '  +1.06%  '

$ws.Range("E23").Value = '  -3.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.41'
$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.13'
$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("D29").Value = '3.970.41'
$ws.Range("E29").Value = '  +1.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.79'
$ws.Range("E30").Value = '  +0.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.43'
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.22'
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.83'
$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.24'
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.94'
$ws.Range("E43").Value = '  -2.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.74'
$ws.Range("E44").Value = '  -1.12%  '

$ws.Range("E45").Value = '  +0.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.14'
$ws.Range("E46").Value = '  +6.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.60'
$ws.Range("E47").Value = '  +1.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.39'
$ws.Range("E48").Value = '  +12.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.37'
$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("E50").Value = '  +2.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '391.84'
$ws.Range("E51").Value = '  -0.13%  '
